$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (2023-10-04 -> 2023-10-05, serial 45203 -> 45204) for every
# data row (rows 2 through 89).
for ($row = 2; $row -le 89; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
